$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Day" labels and dates for rows 23-28 (Day 22 .. Day 27)
$days  = @("Day 22", "Day 23", "Day 24", "Day 25", "Day 26", "Day 27")
$dates = @(45824, 45825, 45826, 45827, 45828, 45829)

for ($i = 0; $i -lt $days.Length; $i++) {
    $row = 23 + $i

    $ws.Range("A$row").Value = $days[$i]
    $ws.Range("B$row").Value = $dates[$i]

    # Copy the date-number formatting from the row above (B22 already
    # carries the existing date style) so the new date cells reuse the
    # same cell style instead of creating a brand-new one.
    $ws.Range("B22").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

$ws.Range("C28").Select()
